$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2-9 (A: numeric id, B: country text, C: result text)
$data = @(
    @(3,  "['Ukraine']",                                                      "['']"),
    @(23, "['Bryansk, Ukraine, Ukrainian']",                                  "['']"),
    @(77, "['Russia']",                                                       "['']"),
    @(84, "['Bakhmut, Donetsk, Pavlohrad, Russia, Ukraine, Ukrainian']",      "['']"),
    @(71, "['Pavlohrad, Ukraine']",                                           "['']"),
    @(39, "['Russia, Ukraine']",                                              "['']"),
    @(73, "['Kyiv, Murmansk, Pavlohrad, Russia, Ukraine, Ukrainian']",        "['']"),
    @(96, "['Russia, Ukraine']",                                              "['']")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 1).Value = $data[$i][0]
    $ws.Cells.Item($r, 2).Value = $data[$i][1]
    $ws.Cells.Item($r, 3).Value = $data[$i][2]
}

# Remove now-unused rows 10-15 that belonged to the old, larger data set
# (-4162 = xlShiftUp)
$ws.Range("A10:C15").Delete(-4162)
